$d = $word.ActiveDocument
$t = $d.Tables(1)

function Replace-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

Replace-CellText 3 2 "Парковка на высоте" "`"НЕ`" Место под колеса"
Replace-CellText 3 31 "Парковка на высоте" "`"НЕ`" Место под колеса"
Replace-CellText 4 2 "ИНН: 624436165689" "ИНН: 701454590182"
Replace-CellText 4 15 "Тверская область, город Солнечногорск, наб. Сталина, 42" "Тверская область, город Видное, бульвар Балканская, 55"
Replace-CellText 4 31 "ИНН: 624436165689" "ИНН: 701454590182"
Replace-CellText 5 2 "Тверская область, город Солнечногорск, наб. Сталина, 42" "Тверская область, город Видное, бульвар Балканская, 55"
Replace-CellText 5 15 "Место №1" "Место №9"
Replace-CellText 5 31 "Тверская область, город Солнечногорск, наб. Сталина, 42" "Тверская область, город Видное, бульвар Балканская, 55"
Replace-CellText 6 6 "000001" "000002"
Replace-CellText 7 14 "Citroen" "Morgan"
Replace-CellText 7 40 "000001" "000002"
Replace-CellText 8 2 "C3" "Aero 8"
Replace-CellText 8 12 "К 474 НК 725" "О 315 ХХ 04"
Replace-CellText 9 3 "Яковлева Арина Ивановна" "Чумакова Анна Арсентьевна"
Replace-CellText 9 10 "г. Москва, ул. Родниковая, 47, оф. 96" "г. Пермь, ул. Луговая, 24, оф. 54"
Replace-CellText 11 3 "79395080159" "79837842573"
Replace-CellText 12 3 "27.01.2025" "28.01.2025"
Replace-CellText 12 8 "21" "15"
Replace-CellText 12 12 "43" "44"
Replace-CellText 13 3 "27.01.2025" "28.01.2025"
Replace-CellText 13 9 "27.02.2025" "07.02.2025"
Replace-CellText 17 7 "Кривоносов Иван Алексеевич" "Чернова Анна Михайловна"
Replace-CellText 19 4 "31" "10"
Replace-CellText 20 2 "6603" "5254,40"

Write-Host "Done"
